$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exigences")

# The "Chapitre" column header (B1) is renamed.
$ws.Range("B1").Value = "ChapitreFFFFFFFFF"

# Move the UI selection from wherever it was (e.g. M22) back to B1.
$ws.Activate()
$ws.Range("B1").Select()
